$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios")

# --- Remove cached formulas from F3 / F4, keeping the literal values ---
# (3408.1+3572 -> 6980.1 ; 56*300 -> 16800)
$ws.Range("F3").Value = 6980.1
$ws.Range("F4").Value = 16800

# --- Add new scenario row 6 ("CH4_Test"), cloned from row 2's layout/styles ---
$ws.Range("A2:L2").Copy()
$ws.Rows.Item(6).Insert()
$ws.Range("A6:D6").ClearFormats()
$ws.Range("M6").Clear()
$ws.Range("A6").Value = "CH4_Test"
$ws.Range("B6").Value = "CH4_Test"

# --- Add new scenario row 7 ("H2_Test"), cloned from row 2's layout/styles ---
$ws.Range("A2:L2").Copy()
$ws.Rows.Item(7).Insert()
$ws.Range("A7:D7").ClearFormats()
$ws.Range("M7").Clear()
$ws.Range("A7").Value = "H2_Test"
$ws.Range("B7").Value = "H2_Test"

# --- Update the active selection to B6 (matches the saved view state) ---
[void]$ws.Range("B6").Select()
